$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update problem-statement text and approval status per the updated
# "approved problem statements" list.
$ws.Range("D4").Value = "Fourier Series Library"
$ws.Range("G4").Value = "Yes"
$ws.Range("G7").Value = "Yes"

# Move the active selection to D5 (reflecting where the author left off).
$ws.Range("D5").Select()
